# Adds an "Electrode Locations" column (C) derived from the file name in
# column A, then re-sorts the data rows (A2:C60) by electrode location
# (letter prefix, then numeric suffix) from A1 through O15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out how many data rows currently exist below the header (row 1).
$lastRow = 1
while ($ws.Cells.Item($lastRow + 1, 1).Value2 -ne $null -and $ws.Cells.Item($lastRow + 1, 1).Value2 -ne "") {
    $lastRow = $lastRow + 1
}

# Read the existing File Name / Unnormalized P_max pairs, and derive the
# electrode location (e.g. "A11") from the leading letters+digits of the
# file name.
$rows = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fname = $ws.Cells.Item($r, 1).Value2
    $val = $ws.Cells.Item($r, 2).Value2

    $letters = ""
    $digits = ""
    if ($fname -match '^([A-Za-z]+)(\d+)') {
        $letters = $matches[1]
        $digits = $matches[2]
    }
    $loc = "$letters$digits"

    $rows += [PSCustomObject]@{
        FName  = $fname
        Val    = $val
        Letter = $letters
        Num    = [int]$digits
        Loc    = $loc
    }
}

# Sort by electrode location: letter prefix first, then numeric suffix
# (natural/numeric order, e.g. A1, A4, A6 ... A11, A12 ... B5, C1 ...).
$sorted = $rows | Sort-Object Letter, Num

# Write the sorted rows back, including the new Electrode Locations column.
$r = 2
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value = $row.FName
    $ws.Cells.Item($r, 2).Value = $row.Val
    $ws.Cells.Item($r, 3).Value = $row.Loc
    $r = $r + 1
}

# Add the new column header, matching the formatting already used by the
# other header cells (bold, centered, bordered).
$ws.Range("C1").Value = "Electrode Locations"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
